$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2520
$ws.Range("C2").Value = 2485
$ws.Range("D2").Value = 2506
$ws.Range("E2").Value = 2509.5
$ws.Range("F2").Value = 42
$ws.Range("G2").Value = 2495.4

$ws.Range("B3").Value = 387.6
$ws.Range("C3").Value = 382.7
$ws.Range("D3").Value = 384.15
$ws.Range("E3").Value = 384.3
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 383.2

$ws.Range("B4").Value = 1522.95
$ws.Range("C4").Value = 1510.1
$ws.Range("D4").Value = 1519
$ws.Range("E4").Value = 1519.55
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 1515.45

$ws.Range("B5").Value = 7375
$ws.Range("C5").Value = 7303.7
$ws.Range("D5").Value = 7363
$ws.Range("E5").Value = 7363.2
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 7332.9

$ws.Range("B6").Value = 239.15
$ws.Range("C6").Value = 232.4
$ws.Range("D6").Value = 238.5
$ws.Range("E6").Value = 238.5
$ws.Range("F6").Value = 106
$ws.Range("G6").Value = 233.85

$ws.Range("B7").Value = 196.75
$ws.Range("C7").Value = 194.55
$ws.Range("D7").Value = 196.2
$ws.Range("E7").Value = 196.25
$ws.Range("F7").Value = 105
$ws.Range("G7").Value = 194.6

$ws.Range("B8").Value = 274.75
$ws.Range("C8").Value = 257.1
$ws.Range("D8").Value = 273.5
$ws.Range("E8").Value = 274
$ws.Range("F8").Value = 722
$ws.Range("G8").Value = 258

$ws.Range("B9").Value = 532.5
$ws.Range("C9").Value = 514.55
$ws.Range("D9").Value = 531.5
$ws.Range("E9").Value = 530.75
$ws.Range("F9").Value = 66
$ws.Range("G9").Value = 518.55

$ws.Range("B10").Value = 3411.95
$ws.Range("C10").Value = 3375.05
$ws.Range("D10").Value = 3397.4
$ws.Range("E10").Value = 3402.45
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = 3382.4

$ws.Range("B11").Value = 147.15
$ws.Range("C11").Value = 144.15
$ws.Range("D11").Value = 145.25
$ws.Range("E11").Value = 145
$ws.Range("F11").Value = 131
$ws.Range("G11").Value = 144.15

$ws.Range("B12").Value = 1257
$ws.Range("C12").Value = 1230.45
$ws.Range("D12").Value = 1255
$ws.Range("E12").Value = 1254
$ws.Range("F12").Value = 25
$ws.Range("G12").Value = 1233.3

$ws.Range("B13").Value = 1614
$ws.Range("C13").Value = 1587.15
$ws.Range("D13").Value = 1611
$ws.Range("E13").Value = 1610.85
$ws.Range("F13").Value = 171
$ws.Range("G13").Value = 1588.5

$ws.Range("B14").Value = 479.5
$ws.Range("C14").Value = 473.55
$ws.Range("D14").Value = 475.1
$ws.Range("E14").Value = 475.1
$ws.Range("F14").Value = 40
$ws.Range("G14").Value = 476.15

$ws.Range("B15").Value = 969.5
$ws.Range("C15").Value = 955.05
$ws.Range("D15").Value = 962.95
$ws.Range("E15").Value = 965.65
$ws.Range("F15").Value = 134
$ws.Range("G15").Value = 956.7

$ws.Range("B16").Value = 1443.2
$ws.Range("C16").Value = 1401.6
$ws.Range("D16").Value = 1438
$ws.Range("E16").Value = 1439.6
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 1403.25

$ws.Range("B17").Value = 1472
$ws.Range("C17").Value = 1462.05
$ws.Range("D17").Value = 1467.3
$ws.Range("E17").Value = 1466.2
$ws.Range("F17").Value = 44
$ws.Range("G17").Value = 1469.9

$ws.Range("B18").Value = 707.75
$ws.Range("C18").Value = 700.1
$ws.Range("D18").Value = 701.8
$ws.Range("E18").Value = 703.5
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 701.05

$ws.Range("B19").Value = 454.95
$ws.Range("C19").Value = 444.95
$ws.Range("D19").Value = 451
$ws.Range("E19").Value = 451
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 449.85

$ws.Range("B20").Value = 1574.45
$ws.Range("C20").Value = 1560.25
$ws.Range("D20").Value = 1564.9
$ws.Range("E20").Value = 1565.75
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 1571.8

$ws.Range("B21").Value = 301
$ws.Range("C21").Value = 293.3
$ws.Range("D21").Value = 300.9
$ws.Range("E21").Value = 299.9
$ws.Range("F21").Value = 33
$ws.Range("G21").Value = 294

$ws.Range("B22").Value = 2438.25
$ws.Range("C22").Value = 2411
$ws.Range("D22").Value = 2431
$ws.Range("E22").Value = 2432
$ws.Range("F22").Value = 68
$ws.Range("G22").Value = 2419

$ws.Range("B23").Value = 581
$ws.Range("C23").Value = 571.9
$ws.Range("D23").Value = 580.6
$ws.Range("E23").Value = 579.05
$ws.Range("F23").Value = 123
$ws.Range("G23").Value = 572.15

$ws.Range("B24").Value = 623.85
$ws.Range("C24").Value = 612.6
$ws.Range("D24").Value = 618
$ws.Range("E24").Value = 617.35
$ws.Range("F24").Value = 8
$ws.Range("G24").Value = 620.65

$ws.Range("B25").Value = 1088.25
$ws.Range("C25").Value = 1078.2
$ws.Range("D25").Value = 1085.15
$ws.Range("E25").Value = 1086.15
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 1082.35

$ws.Range("B26").Value = 616.5
$ws.Range("C26").Value = 609.05
$ws.Range("D26").Value = 614.25
$ws.Range("E26").Value = 614.9
$ws.Range("F26").Value = 65
$ws.Range("G26").Value = 609.2

$ws.Range("B27").Value = 266.35
$ws.Range("C27").Value = 256.3
$ws.Range("D27").Value = 263.65
$ws.Range("E27").Value = 263.7
$ws.Range("F27").Value = 177
$ws.Range("G27").Value = 257.2

$ws.Range("B28").Value = 130.8
$ws.Range("C28").Value = 128.75
$ws.Range("D28").Value = 130.1
$ws.Range("E28").Value = 130.15
$ws.Range("F28").Value = 302
$ws.Range("G28").Value = 128.8

$ws.Range("B29").Value = 8515.75
$ws.Range("C29").Value = 8430
$ws.Range("D29").Value = 8482.5
$ws.Range("E29").Value = 8495.15
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 8475.4

